$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value (kept as text to match original inline-string formatting)
$updates = @{
    'D2' = '304.51'
    'E2' = '2.12%'
    'D3' = '31.87'
    'E3' = '-0.16%'
    'D4' = '5.233'
    'E4' = '2.83%'
    'D5' = '0.07829'
    'E5' = '3.92%'
    'D6' = '2.273'
    'E6' = '31.62%'
    'D7' = '7.998'
    'E7' = '3.20%'
    'D8' = '3.873'
    'E8' = '1.99%'
    'D9' = '0.9167'
    'E9' = '-1.32%'
    'E10' = '2.74%'
    'D11' = '0.07500'
    'E11' = '0.82%'
    'D12' = '0.08215'
    'E12' = '3.66%'
    'D13' = '0.03044'
    'E13' = '-0.43%'
    'D14' = '0.09966'
    'E14' = '0.75%'
    'D15' = '0.001514'
    'E15' = '1.77%'
    'D16' = '0.006174'
    'E16' = '-4.35%'
    'E17' = '1.07%'
    'D18' = '2.240'
    'E18' = '0.79%'
    'E19' = '-0.50%'
    'D20' = '0.1328'
    'E20' = '0.15%'
    'D21' = '4.659'
    'E21' = '2.25%'
    'D22' = '0.04632'
    'E22' = '-0.52%'
    'E23' = '0.34%'
    'D24' = '0.001263'
    'E24' = '3.46%'
    'D25' = '0.004537'
    'E25' = '2.61%'
    'D26' = '0.0001298'
    'E26' = '-7.27%'
    'D27' = '0.0002737'
    'E27' = '47.32%'
    'D39' = '0.01788'
    'E39' = '6.60%'
    'D40' = '0.04588'
    'E40' = '1.07%'
    'D41' = '0.007279'
    'E41' = '2.43%'
    'E42' = '2.92%'
    'D43' = '0.002237'
    'E43' = '8.66%'
    'D44' = '0.01078'
    'E44' = '-8.01%'
    'D45' = '0.00006506'
    'E45' = '8.73%'
    'E46' = '-57.48%'
    'D47' = '0.009882'
    'E47' = '-23.79%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Leading apostrophe forces text interpretation so numeric-looking
    # strings (e.g. '304.51' or '2.12%') are not coerced into numbers.
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = "Normal"
}
